$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header labels in row 2: the placeholder "unnamed" labels
# should read "total" (matching column C, which already says "total").
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
